$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03018965827057
$ws.Range("D2").Value = 1.033834211653507
$ws.Range("E2").Value = 1.040268634461493
$ws.Range("F2").Value = 1.052117186046802
$ws.Range("I2").Value = 1.036122513893768
$ws.Range("J2").Value = 1.035332276107717
$ws.Range("K2").Value = 1.036635421405562
$ws.Range("L2").Value = 1.043051463757426
$ws.Range("M2").Value = 1.054866806748396
$ws.Range("N2").Value = 1.015781473368811
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030942519852846
$ws.Range("D3").Value = 1.03437577644047
$ws.Range("E3").Value = 1.040973717895764
$ws.Range("F3").Value = 1.052985740421213
$ws.Range("I3").Value = 1.036270697065648
$ws.Range("J3").Value = 1.03572749104819
$ws.Range("K3").Value = 1.03698676251411
$ws.Range("L3").Value = 1.043567213091242
$ws.Range("M3").Value = 1.055547999583379
$ws.Range("N3").Value = 1.015913245635315
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031430385047041
$ws.Range("D4").Value = 1.034726779170171
$ws.Range("E4").Value = 1.041431020714989
$ws.Range("F4").Value = 1.053549121343055
$ws.Range("I4").Value = 1.036365761550907
$ws.Range("J4").Value = 1.035983241863133
$ws.Range("K4").Value = 1.037213967982992
$ws.Range("L4").Value = 1.043901328062478
$ws.Range("M4").Value = 1.055989501454468
$ws.Range("N4").Value = 1.015998492824584
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031635652436869
$ws.Range("D5").Value = 1.034874476286533
$ws.Range("E5").Value = 1.041623524168833
$ws.Range("F5").Value = 1.053786292042105
$ws.Range("I5").Value = 1.03640552962477
$ws.Range("J5").Value = 1.036090763069549
$ws.Range("K5").Value = 1.037309451393527
$ws.Range("L5").Value = 1.044041881963562
$ws.Range("M5").Value = 1.056175280431425
$ws.Range("N5").Value = 1.016034325900484
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031670127580944
$ws.Range("D6").Value = 1.034899283177593
$ws.Range("E6").Value = 1.041655861143537
$ws.Range("F6").Value = 1.053826133061903
$ws.Range("I6").Value = 1.036412195283561
$ws.Range("J6").Value = 1.036108816524367
$ws.Range("K6").Value = 1.037325481460967
$ws.Range("L6").Value = 1.044065486929686
$ws.Range("M6").Value = 1.056206483546902
$ws.Range("N6").Value = 1.016040342134062
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031433127178854
$ws.Range("D7").Value = 1.034728752175935
$ws.Range("E7").Value = 1.041433591962614
$ws.Range("F7").Value = 1.053552289153841
$ws.Range("I7").Value = 1.036366293708996
$ws.Range("J7").Value = 1.035984678554386
$ws.Range("K7").Value = 1.037215243971288
$ws.Range("L7").Value = 1.043903205791173
$ws.Range("M7").Value = 1.055991983170364
$ws.Range("N7").Value = 1.015998971647512
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030443942490991
$ws.Range("D8").Value = 1.034017116013533
$ws.Range("E8").Value = 1.040506698630681
$ws.Range("F8").Value = 1.052410433465123
$ws.Range("I8").Value = 1.036172762417368
$ws.Range("J8").Value = 1.035465835857438
$ws.Range("K8").Value = 1.036754186113021
$ws.Range("L8").Value = 1.043225681832634
$ws.Range("M8").Value = 1.055096867907534
$ws.Range("N8").Value = 1.015826009927543
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028706421968321
$ws.Range("D9").Value = 1.032767601797968
$ws.Range("E9").Value = 1.038881653033016
$ws.Range("F9").Value = 1.050408912936342
$ws.Range("I9").Value = 1.035825488694164
$ws.Range("J9").Value = 1.034551781935382
$ws.Range("K9").Value = 1.035940760255191
$ws.Range("L9").Value = 1.042034859054382
$ws.Range("M9").Value = 1.053525192477457
$ws.Range("N9").Value = 1.015521109915208
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027551921484984
$ws.Range("D10").Value = 1.031937716143216
$ws.Range("E10").Value = 1.037803962994768
$ws.Range("F10").Value = 1.049081806709532
$ws.Range("I10").Value = 1.035589817338837
$ws.Range("J10").Value = 1.033942633877023
$ws.Range("K10").Value = 1.035397891843715
$ws.Range("L10").Value = 1.041243128579632
$ws.Range("M10").Value = 1.052481307543023
$ws.Range("N10").Value = 1.015317790891726
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027052946845907
$ws.Range("D11").Value = 1.031579131199156
$ws.Range("E11").Value = 1.037338682699335
$ws.Range("F11").Value = 1.048508900211152
$ws.Range("I11").Value = 1.035486792830788
$ws.Range("J11").Value = 1.033678935025086
$ws.Range("K11").Value = 1.035162701396633
$ws.Range("L11").Value = 1.040900830619274
$ws.Range("M11").Value = 1.052030242321617
$ws.Range("N11").Value = 1.015229745008738
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02686774741036
$ws.Range("D12").Value = 1.031446053061511
$ws.Range("E12").Value = 1.037166064294895
$ws.Range("F12").Value = 1.048296360795618
$ws.Range("I12").Value = 1.035448378873009
$ws.Range("J12").Value = 1.033580996826524
$ws.Range("K12").Value = 1.035075323750091
$ws.Range("L12").Value = 1.040773766636109
$ws.Range("M12").Value = 1.051862840350978
$ws.Range("N12").Value = 1.015197040214814
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026907466861959
$ws.Range("D13").Value = 1.03147459348089
$ws.Range("E13").Value = 1.03720308210905
$ws.Range("F13").Value = 1.048341939258447
$ws.Range("I13").Value = 1.035456625395142
$ws.Range("J13").Value = 1.033602004369908
$ws.Range("K13").Value = 1.035094067320124
$ws.Range("L13").Value = 1.040801018613224
$ws.Range("M13").Value = 1.051898742106511
$ws.Range("N13").Value = 1.015204055525068
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027037635295311
$ws.Range("D14").Value = 1.031568128535564
$ws.Range("E14").Value = 1.037324409764874
$ws.Range("F14").Value = 1.048491326246993
$ws.Range("I14").Value = 1.035483620496442
$ws.Range("J14").Value = 1.03367083918486
$ws.Range("K14").Value = 1.035155479081344
$ws.Range("L14").Value = 1.040890325814811
$ws.Range("M14").Value = 1.052016401879138
$ws.Range("N14").Value = 1.015227041630406
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027117855234019
$ws.Range("D15").Value = 1.03162577403557
$ws.Range("E15").Value = 1.037399191284988
$ws.Range("F15").Value = 1.048583403493862
$ws.Range("I15").Value = 1.035500233736392
$ws.Range("J15").Value = 1.033713252123083
$ws.Range("K15").Value = 1.035193314622226
$ws.Range("L15").Value = 1.040945361673574
$ws.Range("M15").Value = 1.05208891505263
$ws.Range("N15").Value = 1.015241204059033
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027585056587596
$ws.Range("D16").Value = 1.03196153044467
$ws.Range("E16").Value = 1.037834871087439
$ws.Range("F16").Value = 1.049119865464296
$ws.Range("I16").Value = 1.035596634210744
$ws.Range("J16").Value = 1.033960136212038
$ws.Range("K16").Value = 1.035413498140967
$ws.Range("L16").Value = 1.041265856999431
$ws.Range("M16").Value = 1.052511263307042
$ws.Range("N16").Value = 1.015323634091967
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02787837049485
$ws.Range("D17").Value = 1.032172346688588
$ws.Range("E17").Value = 1.038108528876631
$ws.Range("F17").Value = 1.049456841130794
$ws.Range("I17").Value = 1.035656842573608
$ws.Range("J17").Value = 1.034115018768069
$ws.Range("K17").Value = 1.035551580948828
$ws.Range("L17").Value = 1.041467037362141
$ws.Range("M17").Value = 1.052776445268029
$ws.Range("N17").Value = 1.015375338630345
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028049545295629
$ws.Range("D18").Value = 1.032295385534433
$ws.Range("E18").Value = 1.038268280551972
$ws.Range("F18").Value = 1.049653560997012
$ws.Range("I18").Value = 1.035691866751391
$ws.Range("J18").Value = 1.034205365447628
$ws.Range("K18").Value = 1.035632110180956
$ws.Range("L18").Value = 1.04158443316144
$ws.Range("M18").Value = 1.052931212356848
$ws.Range("N18").Value = 1.015405496266286
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028107926643774
$ws.Range("D19").Value = 1.032337350974642
$ws.Range("E19").Value = 1.03832277402253
$ws.Range("F19").Value = 1.049720665757085
$ws.Range("I19").Value = 1.035703793066032
$ws.Range("J19").Value = 1.034236172359288
$ws.Range("K19").Value = 1.035659566470015
$ws.Range("L19").Value = 1.041624470645378
$ws.Range("M19").Value = 1.052983999300325
$ws.Range("N19").Value = 1.015415779099575
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027846891389761
$ws.Range("D20").Value = 1.032149720505256
$ws.Range("E20").Value = 1.038079154346514
$ws.Range("F20").Value = 1.049420669472596
$ws.Range("I20").Value = 1.035650392538852
$ws.Range("J20").Value = 1.034098400673791
$ws.Range("K20").Value = 1.035536767215906
$ws.Range("L20").Value = 1.041445447357692
$ws.Range("M20").Value = 1.052747984333237
$ws.Range("N20").Value = 1.015369791296463
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026999300001238
$ws.Range("D21").Value = 1.031540581567888
$ws.Range("E21").Value = 1.037288676051776
$ws.Range("F21").Value = 1.048447328208333
$ws.Range("I21").Value = 1.035475675138088
$ws.Range("J21").Value = 1.033650568724939
$ws.Range("K21").Value = 1.035137395319272
$ws.Range("L21").Value = 1.040864024828667
$ws.Range("M21").Value = 1.051981750025331
$ws.Range("N21").Value = 1.015220272808915
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026467207988617
$ws.Range("D22").Value = 1.031158265439597
$ws.Range("E22").Value = 1.036792872043075
$ws.Range("F22").Value = 1.047836875959361
$ws.Range("I22").Value = 1.03536497864136
$ws.Range("J22").Value = 1.033369065023033
$ws.Range("K22").Value = 1.034886194088424
$ws.Range("L22").Value = 1.040498929226331
$ws.Range("M22").Value = 1.051500820675529
$ws.Range("N22").Value = 1.01512626117756
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02674920142054
$ws.Range("D23").Value = 1.031360873946846
$ws.Range("E23").Value = 1.037055592581847
$ws.Range("F23").Value = 1.048160342840354
$ws.Range("I23").Value = 1.035423740767666
$ws.Range("J23").Value = 1.03351828875696
$ws.Range("K23").Value = 1.035019369679767
$ws.Range("L23").Value = 1.040692428431691
$ws.Range("M23").Value = 1.051755690859246
$ws.Range("N23").Value = 1.015176098686644
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027861115165267
$ws.Range("D24").Value = 1.032159944077361
$ws.Range("E24").Value = 1.038092427025565
$ws.Range("F24").Value = 1.049437013370204
$ws.Range("I24").Value = 1.035653307323521
$ws.Range("J24").Value = 1.034105909655977
$ws.Range("K24").Value = 1.035543460941832
$ws.Range("L24").Value = 1.041455202794942
$ws.Range("M24").Value = 1.052760844323861
$ws.Range("N24").Value = 1.015372297900233
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029154943156983
$ws.Range("D25").Value = 1.033090088762127
$ws.Range("E25").Value = 1.039300775088976
$ws.Range("F25").Value = 1.050925087302167
$ws.Range("I25").Value = 1.035916002981669
$ws.Range("J25").Value = 1.034788053492107
$ws.Range("K25").Value = 1.036151158951619
$ws.Range("L25").Value = 1.04234234292266
$ws.Range("M25").Value = 1.05393082928297
$ws.Range("N25").Value = 1.015599945144781
